# DEI-3-4 Implemented generic writer
#
# Row 10 (task DEI-3-2, "Clean and transform") gets its real-hours estimate
# bumped from "2 hour 5 min" to "2 hour 40 min" (column B), and the old
# "40 min" note that had been sitting in the Order column (E10) is cleared
# out now that the writer folds it into the single B10 estimate.
#
# Finally, move the cursor down to where the author ended up (D16) after
# making the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Roll the "2 hour 5 min" estimate forward into "2 hour 40 min".
$ws.Range("B10").Value = "2 hour 40 min"

# The separate "40 min" note in the Order column is no longer needed now
# that it has been merged into the estimate above - clear it but keep the
# cell's existing formatting/style intact.
$ws.Range("E10").ClearContents()

# Leave the selection where the author left it.
$ws.Range("D16").Select()
